# Consider `u` for `Jumper`.
#
# The jumper J9-10 (bus 9 <-> bus 10) is converted into a regular (but
# disabled, u=0) Line record "Line_12" inserted at the top of the Line
# table, pushing the existing Line_13..Line_20 rows down by one. Final
# selection/active-sheet state moves from the Jumper sheet to the Line
# sheet.

$wb = $excel.ActiveWorkbook
$lineWs = $wb.Worksheets.Item("Line")
$jumperWs = $wb.Worksheets.Item("Jumper")

# Insert a new row 13 on the Line sheet, shifting the existing Line_13..
# Line_20 rows down to 14..21.
$lineWs.Rows.Item(13).Insert()

# Carry over the formatting (bold/bordered "uid" style) from the cell
# that used to be row 13 (now row 14) onto the freshly inserted row.
$lineWs.Range("A14").Copy()
$lineWs.Range("A13").PasteSpecial(-4122)

# Populate the new row with the former jumper, now expressed as a
# (disabled) line between bus 9 and bus 10.
$lineWs.Range("A13").Value = 11
$lineWs.Range("B13").Value = "Line_12"
$lineWs.Range("C13").Value = 0
$lineWs.Range("D13").Value = "Line_12"
$lineWs.Range("E13").Value = 9
$lineWs.Range("F13").Value = 10
$lineWs.Range("G13").Value = 100
$lineWs.Range("H13").Value = 60
$lineWs.Range("I13").Value = 138
$lineWs.Range("J13").Value = 138
$lineWs.Range("K13").Value = 0.031809999999999998
$lineWs.Range("L13").Value = 0.084500000000000006
$lineWs.Range("M13").Value = 0
$lineWs.Range("N13").Value = 0
$lineWs.Range("O13").Value = 0
$lineWs.Range("P13").Value = 0
$lineWs.Range("Q13").Value = 0
$lineWs.Range("R13").Value = 0
$lineWs.Range("S13").Value = 0
$lineWs.Range("T13").Value = 1
$lineWs.Range("U13").Value = 0

# Move the active selection: Jumper sheet's selection resets to C3 and
# loses tab focus; Line sheet becomes the active tab with C14 selected.
$jumperWs.Select()
$jumperWs.Range("C3").Select()

$lineWs.Select()
$lineWs.Range("C14").Select()
